$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 currently holds the boolean TRUE, which Excel tables render oddly as a
# header. Replace it with the literal text "truez" instead.
$ws.Range("D1").Value = "truez"

# Move the active selection to D1 to match the saved view state.
$ws.Range("D1").Select()
